$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 76
$ws.Range("A76").Value = 46045
$ws.Range("A76").NumberFormat = "d-mmm-yy"
$ws.Range("B76").Value = 5603
$ws.Range("C76").Value = 4203
$ws.Range("D76").Value = 3876
$ws.Range("E76").Value = 243
$ws.Range("F76").Value = 47
$ws.Range("G76").Value = 27
$ws.Range("H76").Value = 9
$ws.Range("I76").Value = 1

# New row 77
$ws.Range("A77").Value = 46048
$ws.Range("A77").NumberFormat = "d-mmm-yy"
$ws.Range("B77").Value = 5599
$ws.Range("C77").Value = 4320
$ws.Range("D77").Value = 3704
$ws.Range("E77").Value = 398
$ws.Range("F77").Value = 134
$ws.Range("G77").Value = 72
$ws.Range("H77").Value = 9
$ws.Range("I77").Value = 3

# Scroll / selection update as recorded in the saved view
$ws.Application.ActiveWindow.ScrollRow = 61
$ws.Range("K76").Select()
